$wb = $excel.ActiveWorkbook

# Touch "list_diseases" first (cursor ends up on B31 there) so that the
# later activation of "general" is what sticks as the saved active tab.
$wsDiseases = $wb.Worksheets.Item("list_diseases")
$wsDiseases.Activate() | Out-Null
$wsDiseases.Range("B31").Select() | Out-Null

# On "general", select & delete the obsolete "runs / number of bootstrap
# runs" row (row 6) -- its row entirely disappears and everything below
# shifts up one row.
$wsGeneral = $wb.Worksheets.Item("general")
$wsGeneral.Activate() | Out-Null
$wsGeneral.Rows("6:6").Select() | Out-Null
$wsGeneral.Rows("6:6").Delete() | Out-Null
$wsGeneral.Range("A6:XFD6").Select() | Out-Null

# Shrink the conditional-formatting ranges that used to reach the
# now-deleted row 6 (the engine doesn't auto-shift these on row delete).
$fcs = $wsGeneral.Range("D2:D11").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($wsGeneral.Range("D2:D5")) | Out-Null
$fcs.Item(2).ModifyAppliesToRange($wsGeneral.Range("D2:D8")) | Out-Null
$fcs.Item(3).ModifyAppliesToRange($wsGeneral.Range("D2:D10")) | Out-Null
$fcs.Item(4).ModifyAppliesToRange($wsGeneral.Range("D2:D10")) | Out-Null
